$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column AW ("dateTNRS") holds a date serial value (45901 = 2025-09-01)
# that must be updated to 45905 (2025-09-05) for all data rows (2-18).
$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(45905)

for ($row = 2; $row -le 18; $row++) {
    $ws.Range("AW$row").Value = $newDate
}
